$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 288
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 97
$ws.Range("N2").Value = 66
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 24
$ws.Range("T2").Value = 43
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 504
$ws.Range("X2").Value = 468
$ws.Range("Z2").Value = 5
$ws.Range("AA2").Value = 5
